# Swap the contents of columns B:AC between paired rows (keep column A intact).
# These pairs of rows had their match data swapped (likely a sort/order fix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(123,124),
    @(180,182),
    @(187,188),
    @(189,190),
    @(195,196),
    @(199,200),
    @(205,206),
    @(211,212),
    @(217,218),
    @(230,231)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}
